# Update the wrapper element used when assembling the MODS update XML.
#
# Row 2 of Sheet1 holds literal XML fragments that get concatenated across
# the row to build one MODS update record per row further down. The
# opening wrapper fragment (column C) switches from the old
#   "><update type="MODS"><mods:mods ...>
# to the new
#   "><datastream type="md_descriptive" operation="update"><mods:mods ...>
# and the matching closing wrapper fragment (column V) is updated to match:
#   </mods:mods></update></object>  ->  </mods:mods></datastream></object>

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">'
$ws.Range("V2").Value = '</mods:mods></datastream></object>'

# Match the author's final selection on the sheet (cell V2, the cell that
# was just edited).
$ws.Activate()
$ws.Range("V2").Select()
